$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (GO:0030047 - actin modification)
$ws.Range("A2").Value = "GO:0030047"
$ws.Range("B2").Value = "actin modification"
$ws.Range("C2").Value = 0.000453
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 46.86
$ws.Range("F2").Value = 7404
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 158
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = "[TRIM32  -  tripartite motif containing 32, TRIOBP  -  trio and f-actin binding protein]"

# Row 3 (GO:1903506 - regulation of nucleic acid-templated transcription)
$ws.Range("A3").Value = "GO:1903506"
$ws.Range("B3").Value = "regulation of nucleic acid-templated transcription"
$ws.Range("C3").Value = 0.000585
$ws.Range("C3").NumberFormat = "0.00E+00"
$ws.Range("D3").Value = 1
$ws.Range("D3").NumberFormat = "0.00E+00"
$ws.Range("E3").Value = 1.56
$ws.Range("F3").Value = 7404
$ws.Range("G3").Value = 1439
$ws.Range("H3").Value = 158
$ws.Range("I3").Value = 48
$ws.Range("J3").Value = "[PSMD11  -  proteasome (prosome, macropain) 26s subunit, non-atpase, 11, TRAK1  -  trafficking protein, kinesin binding 1, CIR1  -  corepressor interacting with rbpj, 1, TCF12  -  transcription factor 12, ID1  -  inhibitor of dna binding 1, dominant negative helix-loop-helix protein, SFRP5  -  secreted frizzled-related protein 5, ZBTB7A  -  zinc finger and btb domain containing 7a, NLRC5  -  nlr family, card domain containing 5, NFKBIE  -  nuclear factor of kappa light polypeptide gene enhancer in b-cells inhibitor, epsilon, RUNX1  -  runt-related transcription factor 1, ZNF423  -  zinc finger protein 423, CBFA2T3  -  core-binding factor, runt domain, alpha subunit 2; translocated to, 3, ACVR1  -  activin a receptor, type i, NLK  -  nemo-like kinase, SPDEF  -  sam pointed domain containing ets transcription factor, GLI1  -  gli family zinc finger 1, RASL11A  -  ras-like, family 11, member a, DMRTA2  -  dmrt-like family a2, FKBP8  -  fk506 binding protein 8, 38kda, PRDM15  -  pr domain containing 15, MEF2B  -  myocyte enhancer factor 2b, HAS3  -  hyaluronan synthase 3, DOT1L  -  dot1-like histone h3k79 methyltransferase, IRF2BPL  -  interferon regulatory factor 2 binding protein-like, ZBTB39  -  zinc finger and btb domain containing 39, TSC22D3  -  tsc22 domain family, member 3, PLXND1  -  plexin d1, TBR1  -  t-box, brain, 1, TRIM32  -  tripartite motif containing 32, CXXC5  -  cxxc finger protein 5, SCMH1  -  sex comb on midleg homolog 1 (drosophila), ATXN1L  -  ataxin 1-like, BAZ1A  -  bromodomain adjacent to zinc finger domain, 1a, RNF41  -  ring finger protein 41, GATA2  -  gata binding protein 2, PPHLN1  -  periphilin 1, GAL  -  galanin/gmap prepropeptide, UBTF  -  upstream binding transcription factor, rna polymerase i, SIX4  -  six homeobox 4, NCOR2  -  nuclear receptor corepressor 2, NFIA  -  nuclear factor i/a, ELP3  -  elongator acetyltransferase complex subunit 3, DNMT3B  -  dna (cytosine-5-)-methyltransferase 3 beta, BCL11A  -  b-cell cll/lymphoma 11a (zinc finger protein), TCF3  -  transcription factor 3, HMGN1  -  high mobility group nucleosome binding domain 1, IHH  -  indian hedgehog, NFIC  -  nuclear factor i/c (ccaat-binding transcription factor)]"

# Row 4 (GO:0006355 - regulation of transcription, DNA-templated)
$ws.Range("A4").Value = "GO:0006355"
$ws.Range("B4").Value = "regulation of transcription, DNA-templated"
$ws.Range("C4").Value = 0.000585
$ws.Range("C4").NumberFormat = "0.00E+00"
$ws.Range("D4").Value = 1
$ws.Range("D4").NumberFormat = "0.00E+00"
$ws.Range("E4").Value = 1.56
$ws.Range("F4").Value = 7404
$ws.Range("G4").Value = 1439
$ws.Range("H4").Value = 158
$ws.Range("I4").Value = 48
$ws.Range("J4").Value = "[PSMD11  -  proteasome (prosome, macropain) 26s subunit, non-atpase, 11, TRAK1  -  trafficking protein, kinesin binding 1, CIR1  -  corepressor interacting with rbpj, 1, TCF12  -  transcription factor 12, ID1  -  inhibitor of dna binding 1, dominant negative helix-loop-helix protein, SFRP5  -  secreted frizzled-related protein 5, ZBTB7A  -  zinc finger and btb domain containing 7a, NLRC5  -  nlr family, card domain containing 5, NFKBIE  -  nuclear factor of kappa light polypeptide gene enhancer in b-cells inhibitor, epsilon, RUNX1  -  runt-related transcription factor 1, ZNF423  -  zinc finger protein 423, CBFA2T3  -  core-binding factor, runt domain, alpha subunit 2; translocated to, 3, ACVR1  -  activin a receptor, type i, NLK  -  nemo-like kinase, SPDEF  -  sam pointed domain containing ets transcription factor, GLI1  -  gli family zinc finger 1, RASL11A  -  ras-like, family 11, member a, DMRTA2  -  dmrt-like family a2, FKBP8  -  fk506 binding protein 8, 38kda, PRDM15  -  pr domain containing 15, MEF2B  -  myocyte enhancer factor 2b, HAS3  -  hyaluronan synthase 3, DOT1L  -  dot1-like histone h3k79 methyltransferase, IRF2BPL  -  interferon regulatory factor 2 binding protein-like, ZBTB39  -  zinc finger and btb domain containing 39, TSC22D3  -  tsc22 domain family, member 3, PLXND1  -  plexin d1, TBR1  -  t-box, brain, 1, TRIM32  -  tripartite motif containing 32, CXXC5  -  cxxc finger protein 5, SCMH1  -  sex comb on midleg homolog 1 (drosophila), ATXN1L  -  ataxin 1-like, BAZ1A  -  bromodomain adjacent to zinc finger domain, 1a, RNF41  -  ring finger protein 41, GATA2  -  gata binding protein 2, PPHLN1  -  periphilin 1, GAL  -  galanin/gmap prepropeptide, UBTF  -  upstream binding transcription factor, rna polymerase i, SIX4  -  six homeobox 4, NCOR2  -  nuclear receptor corepressor 2, NFIA  -  nuclear factor i/a, ELP3  -  elongator acetyltransferase complex subunit 3, DNMT3B  -  dna (cytosine-5-)-methyltransferase 3 beta, BCL11A  -  b-cell cll/lymphoma 11a (zinc finger protein), TCF3  -  transcription factor 3, HMGN1  -  high mobility group nucleosome binding domain 1, IHH  -  indian hedgehog, NFIC  -  nuclear factor i/c (ccaat-binding transcription factor)]"

# Row 5 (GO:2001141 - regulation of RNA biosynthetic process)
$ws.Range("A5").Value = "GO:2001141"
$ws.Range("B5").Value = "regulation of RNA biosynthetic process"
$ws.Range("C5").Value = 0.000615
$ws.Range("C5").NumberFormat = "0.00E+00"
$ws.Range("D5").Value = 1
$ws.Range("D5").NumberFormat = "0.00E+00"
$ws.Range("E5").Value = 1.56
$ws.Range("F5").Value = 7404
$ws.Range("G5").Value = 1442
$ws.Range("H5").Value = 158
$ws.Range("I5").Value = 48
$ws.Range("J5").Value = "[PSMD11  -  proteasome (prosome, macropain) 26s subunit, non-atpase, 11, TRAK1  -  trafficking protein, kinesin binding 1, CIR1  -  corepressor interacting with rbpj, 1, TCF12  -  transcription factor 12, ID1  -  inhibitor of dna binding 1, dominant negative helix-loop-helix protein, SFRP5  -  secreted frizzled-related protein 5, ZBTB7A  -  zinc finger and btb domain containing 7a, NLRC5  -  nlr family, card domain containing 5, NFKBIE  -  nuclear factor of kappa light polypeptide gene enhancer in b-cells inhibitor, epsilon, RUNX1  -  runt-related transcription factor 1, ZNF423  -  zinc finger protein 423, CBFA2T3  -  core-binding factor, runt domain, alpha subunit 2; translocated to, 3, ACVR1  -  activin a receptor, type i, NLK  -  nemo-like kinase, SPDEF  -  sam pointed domain containing ets transcription factor, GLI1  -  gli family zinc finger 1, RASL11A  -  ras-like, family 11, member a, DMRTA2  -  dmrt-like family a2, FKBP8  -  fk506 binding protein 8, 38kda, PRDM15  -  pr domain containing 15, MEF2B  -  myocyte enhancer factor 2b, HAS3  -  hyaluronan synthase 3, DOT1L  -  dot1-like histone h3k79 methyltransferase, IRF2BPL  -  interferon regulatory factor 2 binding protein-like, ZBTB39  -  zinc finger and btb domain containing 39, TSC22D3  -  tsc22 domain family, member 3, PLXND1  -  plexin d1, TBR1  -  t-box, brain, 1, TRIM32  -  tripartite motif containing 32, CXXC5  -  cxxc finger protein 5, SCMH1  -  sex comb on midleg homolog 1 (drosophila), ATXN1L  -  ataxin 1-like, BAZ1A  -  bromodomain adjacent to zinc finger domain, 1a, RNF41  -  ring finger protein 41, GATA2  -  gata binding protein 2, PPHLN1  -  periphilin 1, GAL  -  galanin/gmap prepropeptide, UBTF  -  upstream binding transcription factor, rna polymerase i, SIX4  -  six homeobox 4, NCOR2  -  nuclear receptor corepressor 2, NFIA  -  nuclear factor i/a, ELP3  -  elongator acetyltransferase complex subunit 3, DNMT3B  -  dna (cytosine-5-)-methyltransferase 3 beta, BCL11A  -  b-cell cll/lymphoma 11a (zinc finger protein), TCF3  -  transcription factor 3, HMGN1  -  high mobility group nucleosome binding domain 1, IHH  -  indian hedgehog, NFIC  -  nuclear factor i/c (ccaat-binding transcription factor)]"

# Row 6 (GO:0051252 - regulation of RNA metabolic process)
$ws.Range("A6").Value = "GO:0051252"
$ws.Range("B6").Value = "regulation of RNA metabolic process"
$ws.Range("C6").Value = 0.000808
$ws.Range("C6").NumberFormat = "0.00E+00"
$ws.Range("D6").Value = 1
$ws.Range("D6").NumberFormat = "0.00E+00"
$ws.Range("E6").Value = 1.51
$ws.Range("F6").Value = 7404
$ws.Range("G6").Value = 1580
$ws.Range("H6").Value = 158
$ws.Range("I6").Value = 51
$ws.Range("J6").Value = "[TRAK1  -  trafficking protein, kinesin binding 1, SFRP5  -  secreted frizzled-related protein 5, ID1  -  inhibitor of dna binding 1, dominant negative helix-loop-helix protein, NLRC5  -  nlr family, card domain containing 5, ZNF423  -  zinc finger protein 423, NLK  -  nemo-like kinase, PCBP3  -  poly(rc) binding protein 3, FKBP8  -  fk506 binding protein 8, 38kda, DMRTA2  -  dmrt-like family a2, SRPK1  -  srsf protein kinase 1, PRDM15  -  pr domain containing 15, DOT1L  -  dot1-like histone h3k79 methyltransferase, PLXND1  -  plexin d1, TBR1  -  t-box, brain, 1, TRIM32  -  tripartite motif containing 32, CXXC5  -  cxxc finger protein 5, SCMH1  -  sex comb on midleg homolog 1 (drosophila), BAZ1A  -  bromodomain adjacent to zinc finger domain, 1a, PPHLN1  -  periphilin 1, UBTF  -  upstream binding transcription factor, rna polymerase i, SIX4  -  six homeobox 4, NFIA  -  nuclear factor i/a, ELP3  -  elongator acetyltransferase complex subunit 3, DNMT3B  -  dna (cytosine-5-)-methyltransferase 3 beta, BCL11A  -  b-cell cll/lymphoma 11a (zinc finger protein), TCF3  -  transcription factor 3, HMGN1  -  high mobility group nucleosome binding domain 1, NFIC  -  nuclear factor i/c (ccaat-binding transcription factor), PSMD11  -  proteasome (prosome, macropain) 26s subunit, non-atpase, 11, CIR1  -  corepressor interacting with rbpj, 1, TCF12  -  transcription factor 12, ZBTB7A  -  zinc finger and btb domain containing 7a, NFKBIE  -  nuclear factor of kappa light polypeptide gene enhancer in b-cells inhibitor, epsilon, RUNX1  -  runt-related transcription factor 1, CBFA2T3  -  core-binding factor, runt domain, alpha subunit 2; translocated to, 3, ACVR1  -  activin a receptor, type i, SPDEF  -  sam pointed domain containing ets transcription factor, GLI1  -  gli family zinc finger 1, RASL11A  -  ras-like, family 11, member a, MEF2B  -  myocyte enhancer factor 2b, HAS3  -  hyaluronan synthase 3, IRF2BPL  -  interferon regulatory factor 2 binding protein-like, ZBTB39  -  zinc finger and btb domain containing 39, TSC22D3  -  tsc22 domain family, member 3, ATXN1L  -  ataxin 1-like, RNF41  -  ring finger protein 41, GATA2  -  gata binding protein 2, GAL  -  galanin/gmap prepropeptide, NCOR2  -  nuclear receptor corepressor 2, CPSF7  -  cleavage and polyadenylation specific factor 7, 59kda, IHH  -  indian hedgehog]"
